$d = $word.ActiveDocument

# The table's header row reads: Variable | Overall | Group1 | Group2
# Rename the two group header cells:
#   "Group1" -> "Group0"
#   "Group2" -> "Group1"
# Target the specific table cells directly (row 1, columns 3 and 4) so
# the two renames can't collide with one another.
$table = $d.Tables(1)
$table.Cell(1, 3).Range.Text = "Group0"
$table.Cell(1, 4).Range.Text = "Group1"
